$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1, columns B:I) with the new reference-level
# column labels (pessimistic / class-boundary / optimistic scheme).
$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"

# Move the active selection from C11 to A2.
$ws.Range("A2").Select()
